$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (shifts existing rows 10-14 down to 11-15)
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Formula = "=A9 +1"
$ws.Cells.Item(10, 2).Value = "Apply changes on image of patient face for the use of the surgeon"
$ws.Cells.Item(10, 3).Value = "Must "
$ws.Cells.Item(10, 4).Value = 3
$ws.Cells.Item(10, 5).Value = "To be started "

# Fix the formula reference in the row pushed below the inserted row
$ws.Cells.Item(11, 1).Formula = "=A10 +1"

# Fix text on what is now row 13 ("Get ID 7-8" -> "Get ID 7-9")
$ws.Cells.Item(13, 2).Value = "Get ID 7-9 on the python app"

# Restore column B width (closest achievable value to 59.42578125 given the engine's internal
# pixel-rounding of ColumnWidth) and the active selection (B11)
$ws.Columns.Item(2).ColumnWidth = 58.6
$ws.Cells.Item(11, 2).Select() | Out-Null
